# Auto-generated edit script: updates column F ("想去人数") values
# across all four worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 194
$ws.Cells.Item(4, 6).Value = 329
$ws.Cells.Item(5, 6).Value = 403
$ws.Cells.Item(6, 6).Value = 2182
$ws.Cells.Item(8, 6).Value = 1334
$ws.Cells.Item(9, 6).Value = 3056
$ws.Cells.Item(10, 6).Value = 2042
$ws.Cells.Item(11, 6).Value = 1449
$ws.Cells.Item(12, 6).Value = 1768
$ws.Cells.Item(13, 6).Value = 77
$ws.Cells.Item(14, 6).Value = 695
$ws.Cells.Item(15, 6).Value = 620
$ws.Cells.Item(16, 6).Value = 1233
$ws.Cells.Item(17, 6).Value = 2250
$ws.Cells.Item(18, 6).Value = 1260
$ws.Cells.Item(19, 6).Value = 228
$ws.Cells.Item(20, 6).Value = 2235
$ws.Cells.Item(21, 6).Value = 1863
$ws.Cells.Item(23, 6).Value = 5547
$ws.Cells.Item(24, 6).Value = 1053
$ws.Cells.Item(25, 6).Value = 101
$ws.Cells.Item(27, 6).Value = 1271
$ws.Cells.Item(28, 6).Value = 245
$ws.Cells.Item(29, 6).Value = 1062
$ws.Cells.Item(30, 6).Value = 551
$ws.Cells.Item(31, 6).Value = 122
$ws.Cells.Item(33, 6).Value = 1180
$ws.Cells.Item(35, 6).Value = 3517
$ws.Cells.Item(36, 6).Value = 604
$ws.Cells.Item(39, 6).Value = 92
$ws.Cells.Item(40, 6).Value = 936
$ws.Cells.Item(41, 6).Value = 1230
$ws.Cells.Item(44, 6).Value = 808
$ws.Cells.Item(48, 6).Value = 53
$ws.Cells.Item(49, 6).Value = 28

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 75
$ws.Cells.Item(5, 6).Value = 148658
$ws.Cells.Item(6, 6).Value = 21
$ws.Cells.Item(7, 6).Value = 67
$ws.Cells.Item(10, 6).Value = 122
$ws.Cells.Item(11, 6).Value = 255
$ws.Cells.Item(12, 6).Value = 414
$ws.Cells.Item(16, 6).Value = 329
$ws.Cells.Item(17, 6).Value = 101
$ws.Cells.Item(21, 6).Value = 101
$ws.Cells.Item(24, 6).Value = 3
$ws.Cells.Item(26, 6).Value = 538
$ws.Cells.Item(27, 6).Value = 184
$ws.Cells.Item(28, 6).Value = 331
$ws.Cells.Item(29, 6).Value = 15
$ws.Cells.Item(31, 6).Value = 54
$ws.Cells.Item(32, 6).Value = 54
$ws.Cells.Item(35, 6).Value = 287
$ws.Cells.Item(36, 6).Value = 9
$ws.Cells.Item(40, 6).Value = 21
$ws.Cells.Item(41, 6).Value = 100
$ws.Cells.Item(43, 6).Value = 141
$ws.Cells.Item(46, 6).Value = 3
$ws.Cells.Item(47, 6).Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 3266
$ws.Cells.Item(7, 6).Value = 921
$ws.Cells.Item(8, 6).Value = 1410
$ws.Cells.Item(9, 6).Value = 750
$ws.Cells.Item(10, 6).Value = 349
$ws.Cells.Item(11, 6).Value = 2628
$ws.Cells.Item(12, 6).Value = 188
$ws.Cells.Item(13, 6).Value = 300
$ws.Cells.Item(14, 6).Value = 1009

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 921
$ws.Cells.Item(4, 6).Value = 750
$ws.Cells.Item(5, 6).Value = 194
$ws.Cells.Item(6, 6).Value = 329
$ws.Cells.Item(7, 6).Value = 349
$ws.Cells.Item(8, 6).Value = 2628
$ws.Cells.Item(9, 6).Value = 2182
$ws.Cells.Item(11, 6).Value = 1334
$ws.Cells.Item(12, 6).Value = 3056
$ws.Cells.Item(13, 6).Value = 2042
$ws.Cells.Item(14, 6).Value = 1449
$ws.Cells.Item(16, 6).Value = 1768
$ws.Cells.Item(17, 6).Value = 695
$ws.Cells.Item(18, 6).Value = 620
$ws.Cells.Item(19, 6).Value = 329
$ws.Cells.Item(20, 6).Value = 1233
$ws.Cells.Item(21, 6).Value = 2250
$ws.Cells.Item(22, 6).Value = 188
$ws.Cells.Item(23, 6).Value = 1260
$ws.Cells.Item(24, 6).Value = 228
$ws.Cells.Item(25, 6).Value = 2235
$ws.Cells.Item(26, 6).Value = 5547
$ws.Cells.Item(27, 6).Value = 300
$ws.Cells.Item(28, 6).Value = 1053
$ws.Cells.Item(29, 6).Value = 101
$ws.Cells.Item(30, 6).Value = 1009
$ws.Cells.Item(31, 6).Value = 1272
$ws.Cells.Item(32, 6).Value = 331
$ws.Cells.Item(33, 6).Value = 54
$ws.Cells.Item(34, 6).Value = 1062
$ws.Cells.Item(35, 6).Value = 551
$ws.Cells.Item(36, 6).Value = 122
$ws.Cells.Item(38, 6).Value = 1180
$ws.Cells.Item(39, 6).Value = 3517
$ws.Cells.Item(40, 6).Value = 604
$ws.Cells.Item(41, 6).Value = 287
$ws.Cells.Item(43, 6).Value = 92
$ws.Cells.Item(44, 6).Value = 936
$ws.Cells.Item(45, 6).Value = 1230
$ws.Cells.Item(47, 6).Value = 808
$ws.Cells.Item(49, 6).Value = 141
$ws.Cells.Item(50, 6).Value = 141
$ws.Cells.Item(51, 6).Value = 53

Write-Host "Updated F column values across 展览/演出/本地生活/全部类型 sheets"